$d = $word.ActiveDocument

function Replace-DateText($context, $oldDate, $newDate) {
    # Locate $context (which uniquely identifies the spot in the document)
    # without replacing yet, so we can discover exactly where the match
    # landed ...
    $rng = $d.Content
    $found = $rng.Find.Execute($context, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if (-not $found) {
        throw "context not found: $context"
    }
    # ... then narrow the edit down to just the trailing date text so we
    # never touch characters belonging to a differently formatted run
    # that happens to sit next to the date (e.g. bold dates).
    $dateStart = $rng.End - $oldDate.Length
    $sub = $d.Range($dateStart, $rng.End)
    if ($sub.Text -ne $oldDate) {
        throw "unexpected text at match: [$($sub.Text)] vs [$oldDate]"
    }
    $sub.Text = $newDate
}

# 1. "Defendant appeared in Court for sentencing on June 24, 2022."
Replace-DateText "sentencing on June 24, 2022" "June 24, 2022" "June 26, 2022"

# 2. "...shall pay the fines and costs in full by June 24, 2022."
Replace-DateText "in full by June 24, 2022" "June 24, 2022" "June 26, 2022"

# 3. "...show proof of completion ... on or before August 23, 2022."
Replace-DateText "before August 23, 2022" "August 23, 2022" "August 25, 2022"

# 4. "Defendant's driving license is suspended from June 24, 2022, for a term..."
Replace-DateText "suspended from June 24, 2022" "June 24, 2022" "June 26, 2022"
